$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C80").Value = ":`ndie Ansehenen"
$ws.Range("C116").Value = "`ndie Kellner/innen"
$ws.Range("D116").Value = "`nDer Kellner brachte uns die bestellten Getränke."
$ws.Range("E116").Value = "`nEl mesero nos trajo las bebidas pedidas."
$ws.Range("C117").Value = "`nDie Ärzte"
$ws.Range("D117").Value = "`nDer Arzt wird mir helfen, meine Beschwerden zu lindern."
$ws.Range("E117").Value = "`nEl médico me ayudará a aliviar mis molestias."
$ws.Range("C118").Value = "`nDer Unterrichte sind die Plurale der Unterricht."
$ws.Range("D118").Value = "`nDer Unterricht ist heute sehr interessant."
$ws.Range("E118").Value = "`nLa clase es muy interesante hoy."
$ws.Range("C119").Value = "`nDie Schöpfungen"
$ws.Range("E119").Value = "`nLa educación"
$ws.Range("C120").Value = "`ndie Vergangenheiten"
$ws.Range("D120").Value = "`nDie Vergangenheit hat uns gelehrt, wie wir uns verbessern können."
$ws.Range("E120").Value = "`nLa historia nos ha enseñado cómo podemos mejorar."
$ws.Range("C121").Value = " und seiner Flexion`ndie Verwendungen"
$ws.Range("D121").Value = "`nDie Verwendung dieser Technologie ist sehr einfach."
$ws.Range("E121").Value = "`nEl uso de esta tecnología es muy sencillo."
$ws.Range("C122").Value = "`nDie Weihnachtsgesänge"
$ws.Range("E122").Value = "`nLa canción de Navidad"
$ws.Range("C123").Value = "`nDie Gänge."
$ws.Range("D123").Value = "`nDer Gang in meiner Schule ist sehr weitläufig."
$ws.Range("E123").Value = "`nEl pasillo de mi escuela es muy amplio."
$ws.Range("E124").Value = "`nEl poema"
$ws.Range("C125").Value = "`nDie Ferien."
$ws.Range("E125").Value = " ist sehr interessant`nLa vacación es muy interesante."
$ws.Range("C126").Value = "`nDie Briefe."
$ws.Range("E126").Value = "`nEl carta."
$ws.Range("C127").Value = "`nDie Sterne."
$ws.Range("E127").Value = "`nel estrella"
$ws.Range("C128").Value = "`nDie Würste."
$ws.Range("E128").Value = "`nLa salchicha."
$ws.Range("C129").Value = "`ndie Lieblingsfächer"
$ws.Range("D129").Value = "`nMein Lieblingsfach ist Deutsch."
$ws.Range("E129").Value = "`nMi materia favorita es aleman."
$ws.Range("C130").Value = "`ndie Noten"
$ws.Range("D130").Value = "`nIch bin begeistert, dass ich die Note 1,0 erhalten habe."
$ws.Range("E130").Value = "`nEstoy entusiasmado de haber recibido la calificación de 1.0."
$ws.Range("C131").Value = "`ndie Geschichten"
$ws.Range("D131").Value = "`nDie Geschichte hat uns viel zu erzählen."
$ws.Range("E131").Value = "`nLa historia tiene mucho que decirnos."
$ws.Range("E132").Value = "`nLa arte."
